# fix bug hien thi
# - Switch the active sheet/tab from "Name" to "Project".
# - On the "Project" sheet, add a new value "Vệ sinh" in cell E6 (a new
#   step name alongside XNLR/Cắm LED/Hàn gá/Hàn full) and move the
#   selection there.

$wb = $excel.ActiveWorkbook
$wsProject = $wb.Worksheets.Item("Project")

# Put the new label in place.
$wsProject.Range("E6").Value = "Vệ sinh"

# Make "Project" the active/selected sheet and select the new cell,
# matching the tabSelected / activeTab / selection changes in the diff.
$wsProject.Activate()
$wsProject.Range("E6").Select()
